$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 143; this shifts rows 143-192 down to 144-193
$ws.Rows("143").Insert()

# Populate the new row 143 with the new data record
$ws.Range("A143").Value = 9
$ws.Range("B143").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C143").Value = "Metropolitana"
$ws.Range("D143").Value = 45215
$ws.Range("D143").NumberFormat = $ws.Range("D144").NumberFormat
$ws.Range("E143").Value = 13
$ws.Range("F143").Value = 100112022
$ws.Range("G143").Value = "Arveja Verde"
$ws.Range("H143").Value = "Sin especificar"
$ws.Range("I143").Value = "Primera"
$ws.Range("J143").Value = 52
$ws.Range("K143").Value = 8000
$ws.Range("L143").Value = 8000
$ws.Range("M143").Value = 8000
$ws.Range("N143").Value = "`$/malla 10 kilos"
$ws.Range("O143").Value = "Provincia de Melipilla"
$ws.Range("P143").Value = 800
$ws.Range("Q143").Value = 10
$ws.Range("R143").Value = "Hortaliza"
